$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "A1 font name:" $ws.Range("A1").Font.Name
Write-Host "A1 font size:" $ws.Range("A1").Font.Size
Write-Host "A1 font bold:" $ws.Range("A1").Font.Bold
Write-Host "A1 font color:" $ws.Range("A1").Font.Color
Write-Host "A1 HorizontalAlignment:" $ws.Range("A1").HorizontalAlignment
Write-Host "A1 Interior color:" $ws.Range("A1").Interior.Color
Write-Host "A1 Interior ColorIndex:" $ws.Range("A1").Interior.ColorIndex
Write-Host "A1 Interior Pattern:" $ws.Range("A1").Interior.Pattern
Write-Host "---"
Write-Host "A2 borders around:"
Write-Host "A2 Borders(xlEdgeLeft).LineStyle" $ws.Range("A2").Borders.Item(7).LineStyle
